$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Row 18: "billion 2020 dollars" -> "billion 2021 dollars"
$ws.Range("A18").Value = "billion 2021 dollars"

# Row 21: "million 2020 dollars" -> "million 2021 dollars"
$ws.Range("A21").Value = "million 2021 dollars"

# Row 24: "2020 dollars" -> "2021 dollars"
$ws.Range("A24").Value = "2021 dollars"

# Row 26: A26 becomes a formula =1/1.21 (instead of hardcoded constant)
$ws.Range("A26").Formula = "=1/1.21"

# Row 26: B26 "2019 dollars per 2012 dollar" -> "2012 dollars per 2021 dollars"
$ws.Range("B26").Value = "2012 dollars per 2021 dollars"

# Row 29: 'which in this case is "2012 dollars per 2020 dollar."' -> 'which in this case is "2012 dollars per 2021 dollar."'
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2021 dollar."'

# Row 30: stays the same text, just shared-string index changed in the source file
$ws.Range("B30").Value = "2012 dollars are worth more than 2020 dollars, so we need a"

# Update the selection shown when the sheet is active (cosmetic, matches diff)
$ws.Activate()
$ws.Range("A26").Select()
